$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Make room in the sheet for the new entries by inserting rows.
#    (row numbers below refer to the sheet's ORIGINAL layout, before
#    any of these inserts happen)
# ------------------------------------------------------------------

# Insert 9 rows right after the Wednesday column-header row (orig row 14)
# so the Wednesday block gets 2 new data rows while the rest of the
# week keeps its original spacing, just shifted down.
$ws.Rows("15:23").Insert()

# Insert 4 rows right after the existing Tuesday data row (orig row 11)
# for three new Tuesday entries plus a totals row.
$ws.Rows("12:15").Insert()

# Tidy up the blank filler left behind by the inserts (the genuinely
# empty rows/cells shouldn't carry any left-over style/content).
$ws.Range("A15:D15").Clear()
$ws.Range("A21:D27").Clear()

# ------------------------------------------------------------------
# 2) Monday: shift the three existing "Setup" entries 12 hours later
#    (AM -> PM), durations (column D) stay the same.
# ------------------------------------------------------------------
$ws.Range("A3").Value = 0.64583333333333337
$ws.Range("B3").Value = 0.6875
$ws.Range("A4").Value = 0.76041666666666663
$ws.Range("B4").Value = 0.79166666666666663
$ws.Range("A5").Value = 0.91666666666666663
$ws.Range("B5").Value = 0.97916666666666663

# ------------------------------------------------------------------
# 3) Tuesday: three additional entries + a totals row.
# ------------------------------------------------------------------
$ws.Range("A12:B14").Style = "Normal"
$ws.Range("A12").NumberFormat = "h:mm"
$ws.Range("B12").NumberFormat = "h:mm"
$ws.Range("A13").NumberFormat = "h:mm"
$ws.Range("B13").NumberFormat = "h:mm"
$ws.Range("A14").NumberFormat = "h:mm"
$ws.Range("B14").NumberFormat = "h:mm"

$ws.Range("A12").Value = 0.54166666666666663
$ws.Range("B12").Value = 0.57291666666666663
$ws.Range("C12").Value = "Emotiv Research"
$ws.Range("D12").Value = 0.75

$ws.Range("A13").Value = 0.57291666666666663
$ws.Range("B13").Value = 0.625
$ws.Range("C13").Value = "Connecting to Emotiv Headset"
$ws.Range("D13").Value = 1.25

$ws.Range("A14").Value = 0.82291666666666663
$ws.Range("B14").Value = 0.96875
$ws.Range("C14").Value = "Creating Wrapper Classes for Emotiv API"
$ws.Range("D14").Value = 3.5

$ws.Range("D15").Formula = "=SUM(D11:D14)"

# ------------------------------------------------------------------
# 4) Wednesday: two new entries.
# ------------------------------------------------------------------
$ws.Range("A19:B20").Style = "Normal"
$ws.Range("A19").NumberFormat = "h:mm"
$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("A20").NumberFormat = "h:mm"
$ws.Range("B20").NumberFormat = "h:mm"
$ws.Range("C19:D20").Style = "Normal"

$ws.Range("A19").Value = 0.5
$ws.Range("B19").Value = 0.60416666666666663
$ws.Range("C19").Value = "Coding with/Studying Emotiv API"
$ws.Range("D19").Value = 2.5

$ws.Range("A20").Value = 0.60416666666666663
$ws.Range("B20").Value = 0.625
$ws.Range("C20").Value = "Emotiv Research"
$ws.Range("D20").Value = 0.5

# ------------------------------------------------------------------
# 5) Column C needs to be wide enough for the new, longer entries.
# ------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 34.44140625

# ------------------------------------------------------------------
# 6) Update the view: scrolled down a bit, selection on D21.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("D21").Select() | Out-Null

Write-Output "Done"
